$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-69 down to 30-70.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly price record.
$ws.Range("A29").Value = 9
$ws.Range("B29").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44482
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 100112005
$ws.Range("G29").Value = "Puerro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7500
$ws.Range("N29").Value = "$/paquete 20 unidades"
$ws.Range("O29").Value = "Provincia de Chacabuco"
$ws.Range("P29").Value = 375
$ws.Range("Q29").Value = 20
$ws.Range("R29").Value = "Hortaliza"
